$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Formula = "'29.185.01"
$ws.Range("E2").Value = '  -2.46%  '

# Row 3
$ws.Range("D3").Formula = "'1.850.39"

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").Formula = "'0.6954"
$ws.Range("E5").Value = '  -6.05%  '

# Row 6
$ws.Range("D6").Formula = "'238.63"
$ws.Range("E6").Value = '  -1.65%  '

# Row 7
$ws.Range("D7").Formula = "'1.002"
$ws.Range("E7").Value = '  +0.10%  '

# Row 8
$ws.Range("D8").Formula = "'0.3073"
$ws.Range("E8").Value = '  -2.58%  '

# Row 9
$ws.Range("D9").Formula = "'0.07508"
$ws.Range("E9").Value = '  +3.68%  '

# Row 10
$ws.Range("D10").Formula = "'23.49"
$ws.Range("E10").Value = '  -4.71%  '

# Row 11
$ws.Range("D11").Formula = "'0.08116"
$ws.Range("E11").Value = '  -2.66%  '

# Row 12
$ws.Range("D12").Formula = "'1.876.65"
$ws.Range("E12").Value = '  +0.29%  '

# Row 13
$ws.Range("D13").Formula = "'0.7229"
$ws.Range("E13").Value = '  -3.66%  '

# Row 14
$ws.Range("D14").Formula = "'5.173"
$ws.Range("E14").Value = '  -4.04%  '

# Row 15
$ws.Range("D15").Formula = "'89.21"
$ws.Range("E15").Value = '  -3.38%  '

# Row 16
$ws.Range("D16").Formula = "'29.182.93"
$ws.Range("E16").Value = '  -2.49%  '

# Row 17
$ws.Range("D17").Formula = "'5.791"
$ws.Range("E17").Value = '  -5.14%  '

# Row 18
$ws.Range("D18").Formula = "'241.36"
$ws.Range("E18").Value = '  -2.80%  '

# Row 19
$ws.Range("D19").Formula = "'0.000007696"
$ws.Range("E19").Value = '  -1.89%  '

# Row 20
$ws.Range("D20").Formula = "'13.09"
$ws.Range("E20").Value = '  -3.57%  '

# Row 21
$ws.Range("D21").Formula = "'1.001"
$ws.Range("E21").Value = '  -0.05%  '

# Row 22
$ws.Range("D22").Formula = "'2.094.01"
$ws.Range("E22").Value = '  -2.17%  '

# Row 23
$ws.Range("D23").Formula = "'1.002"
$ws.Range("E23").Value = '  +0.18%  '

# Row 24
$ws.Range("D24").Formula = "'7.619"
$ws.Range("E24").Value = '  -4.84%  '

# Row 25
$ws.Range("D25").Formula = "'9.035"
$ws.Range("E25").Value = '  -2.78%  '

# Row 26
$ws.Range("D26").Formula = "'161.51"
$ws.Range("E26").Value = '  -2.22%  '

# Row 27
$ws.Range("D27").Formula = "'0.1456"
$ws.Range("E27").Value = '  -6.67%  '

# Row 28
$ws.Range("E28").Value = '  -3.30%  '

# Row 29
$ws.Range("D29").Formula = "'1.940"
$ws.Range("E29").Value = '  -4.21%  '

# Row 30
$ws.Range("D30").Formula = "'1.392"
$ws.Range("E30").Value = '  -7.14%  '

# Row 31
$ws.Range("D31").Formula = "'1.501"
$ws.Range("E31").Value = '  -2.38%  '

# Row 32
$ws.Range("D32").Formula = "'4.421"
$ws.Range("E32").Value = '  -4.11%  '

# Row 33
$ws.Range("D33").Formula = "'4.043"
$ws.Range("E33").Value = '  -5.16%  '

# Row 34
$ws.Range("D34").Formula = "'0.05234"
$ws.Range("E34").Value = '  -1.81%  '

# Row 35
$ws.Range("D35").Formula = "'1.191"
$ws.Range("E35").Value = '  -3.60%  '

# Row 36
$ws.Range("D36").Formula = "'0.7085"
$ws.Range("E36").Value = '  -5.65%  '

# Row 37
$ws.Range("D37").Formula = "'0.9991"
$ws.Range("E37").Value = '  -0.13%  '

# Row 38
$ws.Range("E38").Value = '  -1.22%  '

# Row 39
$ws.Range("E39").Value = '  -5.26%  '

# Row 40
$ws.Range("D40").Formula = "'2.692"
$ws.Range("E40").Value = '  -2.25%  '

# Row 41
$ws.Range("D41").Formula = "'0.9181"
$ws.Range("E41").Value = '  +6.11%  '

# Row 42
$ws.Range("D42").Formula = "'5.953"
$ws.Range("E42").Value = '  -3.03%  '

# Row 43
$ws.Range("D43").Formula = "'0.4292"
$ws.Range("E43").Value = '  -5.35%  '

# Row 44
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Formula = "'69.93"
$ws.Range("E44").Value = '  -3.33%  '

# Row 45
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Formula = "'1.045.97"
$ws.Range("E45").Value = '  -5.64%  '

# Row 46
$ws.Range("D46").Formula = "'1.002"
$ws.Range("E46").Value = '  -0.01%  '

# Row 47
$ws.Range("D47").Formula = "'102.35"
$ws.Range("E47").Value = '  -2.06%  '

# Row 48
$ws.Range("D48").Formula = "'7.224"
$ws.Range("E48").Value = '  -4.90%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Formula = "'9.298"
$ws.Range("E49").Value = '  -2.55%  '

# Row 50
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Formula = "'1.742"
$ws.Range("E50").Value = '  -6.57%  '

# Row 51
$ws.Range("D51").Formula = "'1.992.67"
$ws.Range("E51").Value = '  -2.33%  '
